$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("D2").Value = 0.0225
$ws.Range("E2").Value = 0.022
$ws.Range("K2").Value = 564.8
$ws.Range("L2").Value = 0.4072979014927525
$ws.Range("M2").Value = 187.1255
$ws.Range("N2").Value = 0.01791771994331457
$ws.Range("O2").Value = 0.3313128541076488
$ws.Range("P2").Value = 187.1255
$ws.Range("Q2").Value = 0.01791771994331457
$ws.Range("R2").Value = 0.3313128541076488
$ws.Range("U2").Value = 15928
$ws.Range("V2").Value = 1.52514458615803
$ws.Range("W2").Value = 0.06794015595269186
$ws.Range("X2").Value = 0.1045500447629817
$ws.Range("Y2").Value = -0.03660988881028984
$ws.Range("Z2").Value = 0.04209890373446754
$ws.Range("AB2").Value = 0.03406246657345421
$ws.Range("AC2").Value = -0.03406246657345421
$ws.Range("AD2").Value = 40068.3
$ws.Range("AF2").Value = 40068.3
$ws.Range("AG2").Value = 24140.3
$ws.Range("AH2").Value = 0.7932447601456291
$ws.Range("AI2").Value = 0.8186543424028636
$ws.Range("AJ2").Value = 0.6980213336263407
$ws.Range("AK2").Value = 0.7311675213002141

# ---- Row 3 (was Luzerner Kantonalbank, becomes Graubuendner Kantonalbank) ----
$ws.Range("B3").Value = "Graubündner Kantonalbank (SWX:GRKP)"
$ws.Range("D3").Value = 0.0225
$ws.Range("E3").Value = 0.0197
$ws.Range("K3").Value = 182.1
$ws.Range("L3").Value = 0.4243766021906315
$ws.Range("M3").Value = 105.078
$ws.Range("N3").Value = 0.02486052949109235
$ws.Range("O3").Value = 0.5770345963756179
$ws.Range("P3").Value = 105.078
$ws.Range("Q3").Value = 0.02486052949109235
$ws.Range("R3").Value = 0.5770345963756179
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 8240.5
$ws.Range("V3").Value = 1.949629734781271
$ws.Range("W3").Value = 0.06794015595269186
$ws.Range("X3").Value = 0.08057685426517866
$ws.Range("Y3").Value = -0.0126366983124868
$ws.Range("Z3").Value = 0.04257240086117092
$ws.Range("AB3").Value = 0.03374348947412184
$ws.Range("AC3").Value = -0.03374348947412184
$ws.Range("AD3").Value = 10674.1
$ws.Range("AF3").Value = 10674.1
$ws.Range("AG3").Value = 2433.6
$ws.Range("AH3").Value = 0.7163440889079782
$ws.Range("AI3").Value = 0.7879774401677223
$ws.Range("AJ3").Value = 0.3653889464438539
$ws.Range("AK3").Value = 0.4586765177073713

# ---- Row 4 (was Graubuendner Kantonalbank, becomes Luzerner Kantonalbank) ----
$ws.Range("B4").Value = "Luzerner Kantonalbank AG (SWX:LUKN)"
$ws.Range("D4").Value = 0.04099999999999999
$ws.Range("E4").Value = 0.0362
$ws.Range("K4").Value = 230.8
$ws.Range("L4").Value = 0.4110418521816563
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("T4").ClearContents()
$ws.Range("W4").Value = 0.08327920906401097
$ws.Range("X4").Value = 0.1344863767999113
$ws.Range("Y4").Value = -0.05120716773590035
$ws.Range("Z4").Value = 0.03208314763390357
$ws.Range("AB4").Value = 0.03406246657345421
$ws.Range("AC4").Value = -0.03406246657345421
$ws.Range("AD4").Value = 20412.8
$ws.Range("AF4").Value = 20412.8
$ws.Range("AG4").Value = 20412.8
$ws.Range("AH4").Value = 0.8418621531558282
$ws.Range("AI4").Value = 0.8628724087788712
$ws.Range("AJ4").Value = 0.8418621531558282
$ws.Range("AK4").Value = 0.8628724087788712

# ---- Row 5 ----
$ws.Range("D5").Value = -0.0229
$ws.Range("E5").Value = 0.022
$ws.Range("K5").Value = 151.9
$ws.Range("L5").Value = 0.3834890179247665
$ws.Range("M5").Value = 82.0475
$ws.Range("N5").Value = 0.03443756558237146
$ws.Range("O5").Value = 0.5401415404871626
$ws.Range("P5").Value = 82.0475
$ws.Range("Q5").Value = 0.03443756558237146
$ws.Range("R5").Value = 0.5401415404871626
$ws.Range("U5").Value = 7687.5
$ws.Range("V5").Value = 3.226652675760755
$ws.Range("W5").Value = 0.0576645660921722
$ws.Range("X5").Value = 0.1045500447629817
$ws.Range("Y5").Value = -0.0468854786708095
$ws.Range("Z5").Value = 0.07392131979695432
$ws.Range("AB5").Value = 0.0348421614283656
$ws.Range("AC5").Value = -0.0348421614283656
$ws.Range("AD5").Value = 8981.4
$ws.Range("AF5").Value = 8981.4
$ws.Range("AG5").Value = 1293.9
$ws.Range("AH5").Value = 0.7903448639991553
$ws.Range("AI5").Value = 0.7649538799601401
$ws.Range("AJ5").Value = 0.3519475573931019
$ws.Range("AK5").Value = 0.3191977501480165

Write-Output "edits applied"
